$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column A width (new column, narrower "date" column got an explicit width)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 11

# ---------------------------------------------------------------------------
# Update existing rows 2-4: the "Schematics + layout / NFET" entry moved up
# to row 2, pushing the two "Silkscreen" entries down to rows 3 and 4. Also
# a new "Implemented?" column (E) value of "yes" is filled in for all of them.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value2 = "Schematics + layout"
$ws.Range("C2").Value2 = "Rewire U5 and U4 as NFET"
$ws.Range("D2").Value2 = "Fixes LED always ON issue, and allows the circuit to work properly"
$ws.Range("E2").Value2 = "yes"
$ws.Range("E2").WrapText = $true

$ws.Range("B3").Value2 = "Silkscreen"
$ws.Range("C3").Value2 = "YQ on silkscreen is backwards, on both front and back side"
$ws.Range("D3").Value2 = "Visual"
$ws.Range("E3").Value2 = "yes"
$ws.Range("E3").WrapText = $true

$ws.Range("B4").Value2 = "Silkscreen"
$ws.Range("C4").Value2 = "Rev 1.0 silkscreen to Rev 2.0 on back"
$ws.Range("D4").Value2 = "visual"
$ws.Range("E4").Value2 = "yes"
$ws.Range("E4").WrapText = $true

$ws.Rows.Item(2).RowHeight = 28.8
$ws.Rows.Item(3).RowHeight = 28.8
$ws.Rows.Item(4).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# New rows 5 and 6: contrast-pin related schematic changes (entered together)
# ---------------------------------------------------------------------------
$ws.Range("A5").Value2 = 42476
$ws.Range("A5").NumberFormat = "d-mmm-yy"
$ws.Range("B5").Value2 = "Schematics"
$ws.Range("B5").NumberFormat = "d-mmm-yy"
$ws.Range("D5").Value2 = "reduces footprint size"
$ws.Range("D5").WrapText = $true

$ws.Range("A6").Value2 = 42476
$ws.Range("A6").NumberFormat = "d-mmm-yy"
$ws.Range("B6").Value2 = "Schematics"
$ws.Range("B6").NumberFormat = "d-mmm-yy"
$ws.Range("C6").Value2 = "Added 10k resistor in parallel to C1 for resistor divider on contrast pin"
$ws.Range("C6").WrapText = $true
$ws.Range("E6").Value2 = "Yes"
$ws.Range("E6").WrapText = $true
$ws.Range("D6").Value2 = "Allow constract pin PWM to actually function correctly"
$ws.Range("D6").WrapText = $true

$ws.Range("C5").Value2 = "Updated caps to use 2.2uF instead of 1000uF, package: 0603"
$ws.Range("C5").WrapText = $true
$ws.Range("E5").Value2 = "yes"
$ws.Range("E5").WrapText = $true

$ws.Rows.Item(5).RowHeight = 28.8
$ws.Rows.Item(6).RowHeight = 28.8

# ---------------------------------------------------------------------------
# New row 7: rotary encoder footprint fix
# ---------------------------------------------------------------------------
$ws.Range("A7").Value2 = 42476
$ws.Range("A7").NumberFormat = "d-mmm-yy"
$ws.Range("B7").Value2 = "Footprint"
$ws.Range("B7").NumberFormat = "d-mmm-yy"
$ws.Range("C7").Value2 = "Increased through-hole pin sizes on rotary encoder"
$ws.Range("C7").WrapText = $true
$ws.Range("D7").Value2 = "Pins actually fit into through-holes on rotary encoder now"
$ws.Range("D7").WrapText = $true
$ws.Range("E7").Value2 = "yes"
$ws.Range("E7").WrapText = $true

$ws.Rows.Item(7).RowHeight = 28.8

# ---------------------------------------------------------------------------
# New row 8: RTC placement tweak
# ---------------------------------------------------------------------------
$ws.Range("A8").Value2 = 42477
$ws.Range("A8").NumberFormat = "d-mmm-yy"
$ws.Range("B8").Value2 = "Schematics + layout"
$ws.Range("B8").NumberFormat = "d-mmm-yy"
$ws.Range("C8").Value2 = "Moved RTC a bit closer to the edge"
$ws.Range("C8").WrapText = $true
$ws.Range("D8").Value2 = "Easier for soldering"
$ws.Range("D8").WrapText = $true
$ws.Range("E8").Value2 = "yes"
$ws.Range("E8").WrapText = $true

$ws.Rows.Item(8).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# New row 9: added testpoints
# ---------------------------------------------------------------------------
$ws.Range("A9").Value2 = 42477
$ws.Range("A9").NumberFormat = "d-mmm-yy"
$ws.Range("B9").Value2 = "Schematics + layout"
$ws.Range("B9").NumberFormat = "d-mmm-yy"
$ws.Range("C9").Value2 = "Addded testpoints"
$ws.Range("C9").WrapText = $true
$ws.Range("D9").Value2 = "for debugging"
$ws.Range("D9").WrapText = $true
$ws.Range("E9").Value2 = "yes"
$ws.Range("E9").WrapText = $true

$ws.Rows.Item(9).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# Freeze the header row and leave the selection on C12, ready for future entry
# ---------------------------------------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Range("C12").Select()

Write-Host "done"
